# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.294.86"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.871.46"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'0.7070"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").Value = "'241.50"
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "'0.07747"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").Value = "'25.08"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").Value = "'0.08379"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "1.864.70"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'5.235"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "'0.7108"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "'91.00"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "29.308.11"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "'6.054"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "'0.000008164"
$ws.Range("E18").Value = "  +4.46%  "
$ws.Range("D19").Value = "'239.51"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").Value = "'13.20"
$ws.Range("D21").Value = "2.118.87"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'7.738"
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'0.1580"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "'163.27"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "'9.009"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").Value = "'18.43"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'1.507"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'1.290"
$ws.Range("D32").Value = "'4.300"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "'0.05288"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").Value = "'0.7434"
$ws.Range("E36").Value = "  -7.51%  "
$ws.Range("D37").Value = "'2.702"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").Value = "'0.01868"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").Value = "1.229.47"
$ws.Range("E39").Value = "  +5.54%  "
$ws.Range("D40").Value = "'2.729"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").Value = "'6.520"
$ws.Range("E41").Value = "  +4.03%  "
$ws.Range("D42").Value = "'0.8846"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'109.27"
$ws.Range("E43").Value = "  +6.34%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'72.27"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D46").Value = "2.016.34"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "'1.792"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").Value = "'0.4302"
$ws.Range("E51").Value = "  +0.24%  "
